$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.693.33'
$ws.Range("D3").Value = '1.644.04'
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = "'212.96"
$ws.Range("E5").Value = '  +0.93%  '
$ws.Range("D6").Value = "'0.529"
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = "'23.39"
$ws.Range("E8").Value = '  +0.97%  '
$ws.Range("E9").Value = '  +1.32%  '
$ws.Range("D11").Value = "'0.0896"
$ws.Range("E11").Value = '  +0.62%  '
$ws.Range("D12").Value = '1.876.71'
$ws.Range("D13").Value = '1.643.93'
$ws.Range("E13").Value = '  +0.32%  '
$ws.Range("D14").Value = "'4.05"
$ws.Range("E14").Value = '  +0.33%  '
$ws.Range("E15").Value = '  +1.22%  '
$ws.Range("D16").Value = "'64.68"
$ws.Range("E16").Value = '  +0.65%  '
$ws.Range("D17").Value = '27.663.01'
$ws.Range("D18").Value = "'230.72"
$ws.Range("E18").Value = '  +0.27%  '
$ws.Range("E20").Value = '  +1.78%  '
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("E22").Value = '  -0.53%  '
$ws.Range("D23").Value = "'10.00"
$ws.Range("E23").Value = '  +7.19%  '
$ws.Range("E24").Value = '  -2.67%  '
$ws.Range("D25").Value = "'149.89"
$ws.Range("E25").Value = '  +1.37%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("E27").Value = '  -1.07%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = "'15.64"
$ws.Range("E28").Value = '  +0.80%  '
$ws.Range("B29").Value = 'BinanceUSD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("E30").Value = '  +0.81%  '
$ws.Range("E31").Value = '  +0.68%  '
$ws.Range("E32").Value = '  +0.73%  '
$ws.Range("D33").Value = '1.443.31'
$ws.Range("E33").Value = '  +2.20%  '
$ws.Range("D34").Value = "'3.13"
$ws.Range("E34").Value = '  +1.02%  '
$ws.Range("E35").Value = '  +1.33%  '
$ws.Range("D36").Value = "'2.35"
$ws.Range("E36").Value = '  -1.04%  '
$ws.Range("E37").Value = '  +1.24%  '
$ws.Range("D38").Value = "'0.880"
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("E39").Value = '  +0.23%  '
$ws.Range("D40").Value = "'0.895"
$ws.Range("E40").Value = '  +13.57%  '
$ws.Range("E41").Value = '  +0.33%  '
$ws.Range("E43").Value = '  +3.43%  '
$ws.Range("D44").Value = "'67.10"
$ws.Range("E44").Value = '  +4.09%  '
$ws.Range("E45").Value = '  +0.37%  '
$ws.Range("D46").Value = '1.786.44'
$ws.Range("E46").Value = '  +0.23%  '
$ws.Range("E47").Value = '  +6.09%  '
$ws.Range("D49").Value = "'0.0987"
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("E50").Value = '  +1.29%  '
$ws.Range("D51").Value = "'0.0503"
$ws.Range("E51").Value = '  +0.70%  '
